$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:E72")
$rng.Sort($ws.Range("E1"), 1, $null, $null, 1, $null, 1, 1)

Write-Host "Done sorting"
